$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")

# Row 17 - add a B17 formula (electricity-adjacent "syngas - ecoinvent" row), mirroring the
# existing C17 formula so B17 now also reports the converted value.
$ws.Range("B17").Formula = "=5.4/1.15"

# New row 20 - charcoal-low upstream
$ws.Range("A20").Value = "charcoal-low upstream"
$ws.Range("C20").Value = 29.5
$ws.Range("D20").Formula = "=112*C20/1000"
$ws.Range("G20").Value = 0.91

# New row 21 - charcoal-high upstream
$ws.Range("A21").Value = "charcoal-high upstream"
$ws.Range("C21").Value = 29.5
$ws.Range("D21").Formula = "=112*C21/1000"
$ws.Range("G21").Value = 0.91

# Move the active selection to D25, matching the saved cursor position.
[void]$ws.Range("D25").Select()
